$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 647 ("「救助」..." entry) entirely; all rows below shift up by one.
$ws.Rows.Item(647).Delete()
